# Apply the "Updated symbol list" edit: refresh crypto price/volume/date/hour
# columns, and shift Coin/Link rows 6-18 by one position (new GateToken entry
# inserted, pushing existing rows down) to match the commit's new snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    @{ Row = 2; D = '313.27'; E = '1.13%'; F = '13-2-2023'; G = '0' },
    @{ Row = 3; D = '40.35'; E = '-2.09%'; F = '13-2-2023'; G = '0' },
    @{ Row = 4; D = '5.139'; E = '-0.29%'; F = '13-2-2023'; G = '0' },
    @{ Row = 5; D = '0.07623'; E = '-0.82%'; F = '13-2-2023'; G = '0' },
    @{ Row = 6; B = 'GateToken'; C = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'; D = '4.333'; E = '0.22%'; F = '13-2-2023'; G = '0' },
    @{ Row = 7; B = 'FTXToken'; C = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; D = '1.713'; E = '5.94%'; F = '13-2-2023'; G = '0' },
    @{ Row = 8; B = 'MXToken'; C = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; D = '0.9390'; E = '1.80%'; F = '13-2-2023'; G = '0' },
    @{ Row = 9; B = 'BTSEToken'; C = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'; D = '2.427'; E = '-2.17%'; F = '13-2-2023'; G = '0' },
    @{ Row = 10; B = 'LiechtensteinCryptoassetsExchange'; C = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'; D = '0.1255'; E = '6.82%'; F = '13-2-2023'; G = '0' },
    @{ Row = 11; B = 'WazirX'; C = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'; D = '0.1821'; E = '-0.98%'; F = '13-2-2023'; G = '0' },
    @{ Row = 12; B = 'MandalaExchangeToken'; C = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D = '0.09014'; E = '-2.14%'; F = '13-2-2023'; G = '0' },
    @{ Row = 13; B = 'BitrueCoin'; C = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D = '0.04130'; E = '-3.92%'; F = '13-2-2023'; G = '0' },
    @{ Row = 14; B = 'BitMartToken'; C = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D = '0.1054'; E = '0.40%'; F = '13-2-2023'; G = '0' },
    @{ Row = 15; B = 'BitForexToken'; C = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D = '0.001268'; E = '0.55%'; F = '13-2-2023'; G = '0' },
    @{ Row = 16; B = 'TigerCash'; C = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D = '0.005866'; E = '0.65%'; F = '13-2-2023'; G = '0' },
    @{ Row = 17; B = 'UpBots'; C = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'; D = '0.007522'; E = '1,905.56%'; F = '13-2-2023'; G = '0' },
    @{ Row = 18; B = 'LEO'; C = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D = '3.359'; F = '13-2-2023'; G = '0' },
    @{ Row = 19; D = '0.3359'; E = '0.70%'; F = '13-2-2023'; G = '0' },
    @{ Row = 20; D = '8.419'; E = '19.65%'; F = '13-2-2023'; G = '0' },
    @{ Row = 21; D = '0.1346'; E = '-3.96%'; F = '13-2-2023'; G = '0' },
    @{ Row = 22; D = '0.2728'; E = '-6.25%'; F = '13-2-2023'; G = '0' },
    @{ Row = 23; D = '0.04030'; E = '-0.69%'; F = '13-2-2023'; G = '0' },
    @{ Row = 24; D = '0.001261'; E = '-0.30%'; F = '13-2-2023'; G = '0' },
    @{ Row = 25; D = '0.004072'; E = '-1.41%'; F = '13-2-2023'; G = '0' },
    @{ Row = 26; D = '0.0001271'; E = '-0.03%'; F = '13-2-2023'; G = '0' },
    @{ Row = 27; F = '13-2-2023'; G = '0' },
    @{ Row = 28; F = '13-2-2023'; G = '0' },
    @{ Row = 29; F = '13-2-2023'; G = '0' },
    @{ Row = 30; F = '13-2-2023'; G = '0' },
    @{ Row = 31; F = '13-2-2023'; G = '0' },
    @{ Row = 32; F = '13-2-2023'; G = '0' },
    @{ Row = 33; F = '13-2-2023'; G = '0' },
    @{ Row = 34; F = '13-2-2023'; G = '0' },
    @{ Row = 35; F = '13-2-2023'; G = '0' },
    @{ Row = 36; F = '13-2-2023'; G = '0' },
    @{ Row = 37; F = '13-2-2023'; G = '0' },
    @{ Row = 38; D = '0.02489'; E = '1.80%'; F = '13-2-2023'; G = '0' },
    @{ Row = 39; D = '0.05209'; E = '-1.40%'; F = '13-2-2023'; G = '0' },
    @{ Row = 40; D = '0.007780'; E = '-0.76%'; F = '13-2-2023'; G = '0' },
    @{ Row = 41; D = '0.1302'; E = '-0.93%'; F = '13-2-2023'; G = '0' },
    @{ Row = 42; D = '0.007678'; E = '12.85%'; F = '13-2-2023'; G = '0' },
    @{ Row = 43; D = '0.002111'; E = '11.02%'; F = '13-2-2023'; G = '0' },
    @{ Row = 44; D = '0.008105'; E = '9.06%'; F = '13-2-2023'; G = '0' },
    @{ Row = 45; D = '0.3127'; E = '-7.29%'; F = '13-2-2023'; G = '0' },
    @{ Row = 46; D = '0.00006614'; E = '-2.77%'; F = '13-2-2023'; G = '0' },
    @{ Row = 47; D = '0.00000000750'; E = '-0.08%'; F = '13-2-2023'; G = '0' },
    @{ Row = 48; D = '0.2682'; E = '30.50%'; F = '13-2-2023'; G = '0' },
    @{ Row = 49; D = '0.004202'; E = '2.49%'; F = '13-2-2023'; G = '0' },
    @{ Row = 50; D = '0.00002101'; E = '-0.08%'; F = '13-2-2023'; G = '0' },
    @{ Row = 51; D = '0.0002001'; E = '-0.08%'; F = '13-2-2023'; G = '0' }
)

foreach ($item in $rowUpdates) {
    $r = $item.Row
    foreach ($col in @("B", "C", "D", "E", "F", "G")) {
        if ($item.ContainsKey($col)) {
            $cellRef = $col + $r
            $ws.Range($cellRef).NumberFormat = "@"
            $ws.Range($cellRef).Value = $item[$col]
        }
    }
}
